# Insert a new data row at row 484 (pushes old rows 484..574 down to 485..575)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(484).Insert()

$ws.Cells.Item(484, 1).Value  = 3
$ws.Cells.Item(484, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(484, 3).Value  = "Coquimbo"
$ws.Cells.Item(484, 4).Value  = 44995
$ws.Cells.Item(484, 5).Value  = 5
$ws.Cells.Item(484, 6).Value  = 100112017
$ws.Cells.Item(484, 7).Value  = "Apio"
$ws.Cells.Item(484, 8).Value  = "Americana (o)"
$ws.Cells.Item(484, 9).Value  = "Primera"
$ws.Cells.Item(484, 10).Value = 240
$ws.Cells.Item(484, 11).Value = 9000
$ws.Cells.Item(484, 12).Value = 9500
$ws.Cells.Item(484, 13).Value = 9229
$ws.Cells.Item(484, 14).Value = "`$/docena de matas"
$ws.Cells.Item(484, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(484, 16).Value = 1538
$ws.Cells.Item(484, 17).Value = 6
$ws.Cells.Item(484, 18).Value = "Hortaliza"
